# Apply the "crypto list" refresh produced by the GitHub Actions job.
# Most rows only get their Price (column D) and/or Volume(1h) (column E)
# values refreshed. Column D is forced to Text format before writing so
# that numeric-looking prices (e.g. "393.61") are stored as exact text,
# matching the workbook's existing inlineStr convention instead of being
# auto-converted to floating point numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.380.33"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.072.30"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "393.61"
$ws.Range("E5").Value = "  +2.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.01"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("E7").Value = "  -1.86%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.11"
$ws.Range("E10").Value = "  +0.46%  "

$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("E12").Value = "  -1.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.558.46"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.36"
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.70"

# Rows 16/17 swap identity: WrappedEther moves up to row 16,
# Polygon moves down to row 17 (with refreshed price/volume figures).
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.071.81"
$ws.Range("E16").Value = "  +1.03%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.02"
$ws.Range("E17").Value = "  +4.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.54"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.388.53"
$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("E20").Value = "  +2.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.30"
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.97"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.83"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  -5.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.95"
$ws.Range("E27").Value = "  +2.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  -4.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.163"
$ws.Range("E30").Value = "  -5.90%  "

$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("E32").Value = "  +4.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0486"
$ws.Range("E33").Value = "  +7.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.17"
$ws.Range("E34").Value = "  +5.95%  "

$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.13"
$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.290"
$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.94"
$ws.Range("E40").Value = "  +7.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.91"
$ws.Range("E41").Value = "  -0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.67"
$ws.Range("E42").Value = "  -1.99%  "

$ws.Range("E43").Value = "  -1.15%  "

$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.75"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("E48").Value = "  -2.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.072.15"
$ws.Range("E49").Value = "  +1.88%  "

$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.894"
$ws.Range("E51").Value = "  +8.31%  "
